# RASD/efforts.xlsx - "add state diagrams and user characteristics"
#
# The third effort-log block (rows 26-33, Sara's log) gets a new combined
# entry for 29/10 + 30/10 + 01/11/2019 covering doc structure, requirements,
# DA, software system attributes, diagrams and user characteristics. The
# hour count for that entry grows from 7 to 9, which ripples into the
# SUM() total in C33 (18 -> 20). The row also grows taller to fit the
# longer wrapped text, and the visible window/selection shifts down a bit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the date list and topic description for row 32 (Sara's last entry)
$ws.Range("A32").Value = "29/10/2019 + 30/10/2019 + 01/11/2019"
$ws.Range("B32").Value = "Doc structure + Requirements + DA + software system attributes + diagrams + user characteristics"

# More topics covered this entry -> more hours logged (7 -> 9). C33 = SUM(C27:C32)
# recalculates automatically to 20.
$ws.Range("C32").Value = 9

# The extra text wraps onto more lines, so the row grows from 58 to 72.5 points.
$ws.Rows.Item(32).RowHeight = 72.5

# Scroll the view down a bit and leave the selection on E30, matching where
# the editor ended up after making the change.
$excel.ActiveWindow.ScrollRow = 22
$ws.Range("E30").Select()
